# Budget creation input data has been updated:
#   - GLCode2 (F2) -> "617392821-RANTY RESER"
#   - GLCode1 (D2) -> "607092921-AED EXPENSES"
# (BudgetReqNumber in C2 keeps displaying "BUD-20FH-0131"; the now-unused
#  "BUD-20FH-0130" string is dropped automatically once nothing references it.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BudgetCreationInputData")
$ws.Activate()

$ws.Range("F2").Value = "617392821-RANTY RESER"
$ws.Range("D2").Value = "607092921-AED EXPENSES"

[void]$ws.Range("D2").Select()
